# ============================================================================
# Edit script: 20201110 results summary both probes.xlsx
# Adds probe C sample data processed on 11/17/2020 (P-0023, B-0041, P-0034,
# BAYSTD1-11172020) across the data sheet and the QAQC assessment sheets.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Workbook-level renames / view state
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data 11dec2019 to 20feb2020")
$wsData.Name = "data 11dec2019 to 17nov2020"

$wsExample   = $wb.Worksheets.Item("example to follow")
$wsCrmBay    = $wb.Worksheets.Item("QAQC crm & baystd assessment")
$wsCrms      = $wb.Worksheets.Item("QAQC crms assessment")
$wsBaystds   = $wb.Worksheets.Item("QAQC baystds assessment")

Write-Output "done"
